$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows 14-17: labels in column A, aggregate formulas in column B ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold / size-12 / vertically-centered look for B14:B17 using a single
# scratch "template" cell so only one combined cell style is produced (rather than
# three separate intermediate styles for Bold, then Size, then VerticalAlignment).
$ws.Range("AA1").Font.Bold = $true
$ws.Range("AA1").Font.Size = 12
$ws.Range("AA1").VerticalAlignment = -4108
$ws.Range("AA1").Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)
$ws.Range("AA1").Clear()

$ws.Range("A14:B17").RowHeight = 15.6

# --- Row 12: bold average of column J (J2:J11) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- Print/page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Match the saved selection ---
[void]$ws.Range("A14:B17").Select()
